$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2785.3845
$ws.Range("J2").Value = 3477.2
$ws.Range("L2").Value = 3477.2
$ws.Range("N2").Value = -3703.2
$ws.Range("H6").Value = 506.33334
$ws.Range("I6").Value = 510
$ws.Range("J6").Value = 499
$ws.Range("K6").Value = 1530
$ws.Range("L6").Value = 1497
$ws.Range("M6").Value = -1418
$ws.Range("N6").Value = -1721
$ws.Range("H8").Value = 136.14285
$ws.Range("I8").Value = 136.14285
$ws.Range("K8").Value = 408.42855
$ws.Range("M8").Value = -269.42855
$ws.Range("H9").Value = 11911043
$ws.Range("J9").Value = 21384.5
$ws.Range("L9").Value = 21384.5
$ws.Range("N9").Value = -21722.5
$ws.Range("H11").Value = 127087.7
$ws.Range("I11").Value = 127087.7
$ws.Range("K11").Value = 127087.7
$ws.Range("M11").Value = -126947.7
$ws.Range("H17").Value = 2117593.8
$ws.Range("J17").Value = 2117593.8
$ws.Range("L17").Value = 6352781.399999999
$ws.Range("N17").Value = -6353117.399999999
$ws.Range("H28").Value = 582.0345
$ws.Range("I28").Value = 611.1539
$ws.Range("K28").Value = 611.1539
$ws.Range("M28").Value = -126.1539
$ws.Range("H32").Value = 23811318
$ws.Range("J32").Value = 33334986
$ws.Range("L32").Value = 33334986
$ws.Range("N32").Value = -33335638
$ws.Range("H33").Value = 23486346
$ws.Range("I33").Value = 14646922
$ws.Range("K33").Value = 14646922
$ws.Range("M33").Value = -14646693
$ws.Range("H40").Value = 2664.5715
$ws.Range("J40").Value = 3237.25
$ws.Range("L40").Value = 3237.25
$ws.Range("N40").Value = -3587.25
$ws.Range("H51").Value = 64103964
$ws.Range("I51").Value = 55557150
$ws.Range("K51").Value = 55557150
$ws.Range("M51").Value = -55556666
$ws.Range("H64").Value = 8428.666999999999
$ws.Range("I64").Value = 7664.1
$ws.Range("K64").Value = 7664.1
$ws.Range("M64").Value = -7416.1
$ws.Range("H67").Value = 8428.666999999999
$ws.Range("I67").Value = 7664.1
$ws.Range("K67").Value = 7664.1
$ws.Range("M67").Value = -6806.1
$ws.Range("H74").Value = 11618.389
$ws.Range("J74").Value = 4332.6665
$ws.Range("L74").Value = 4332.6665
$ws.Range("N74").Value = -6204.6665
$ws.Range("H77").Value = 11618.389
$ws.Range("J77").Value = 4332.6665
$ws.Range("L77").Value = 21663.3325
$ws.Range("N77").Value = -31023.3325
$ws.Range("H86").Value = 83335870
$ws.Range("I86").Value = 83335800
$ws.Range("J86").Value = 83336090
$ws.Range("K86").Value = 83335800
$ws.Range("L86").Value = 83336090
$ws.Range("M86").Value = -83334677
$ws.Range("N86").Value = -83338336
$ws.Range("H89").Value = 83335870
$ws.Range("I89").Value = 83335800
$ws.Range("J89").Value = 83336090
$ws.Range("K89").Value = 416679000
$ws.Range("L89").Value = 416680450
$ws.Range("M89").Value = -416673384
$ws.Range("N89").Value = -416691682
$ws.Range("H101").Value = 246.5
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 72166.86
$ws.Range("I103").Value = 757.8333
$ws.Range("J103").Value = 125723.625
$ws.Range("K103").Value = 2273.4999
$ws.Range("L103").Value = 377170.875
$ws.Range("M103").Value = -1687.4999
$ws.Range("N103").Value = -378342.875
$ws.Range("H133").Value = 126990.2
$ws.Range("J133").Value = 126990.2
$ws.Range("L133").Value = 126990.2
$ws.Range("N133").Value = -137110.2
$ws.Range("H138").Value = 5958.9844
$ws.Range("I138").Value = 15639.091
$ws.Range("J138").Value = 3949.9058
$ws.Range("K138").Value = 46917.273
$ws.Range("L138").Value = 11849.7174
$ws.Range("M138").Value = -41777.273
$ws.Range("N138").Value = -22129.7174
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 223164.45
$ws.Range("I32").Value = 241303.72
$ws.Range("K32").Value = 241303.72
$ws.Range("M32").Value = -241016.72
$ws.Range("H37").Value = 42996.8
$ws.Range("J37").Value = 42996.8
$ws.Range("L37").Value = 42996.8
$ws.Range("N37").Value = -43542.8
$ws.Range("H45").Value = 45960.082
$ws.Range("I45").Value = 59186.39
$ws.Range("J45").Value = 6281.1665
$ws.Range("K45").Value = 59186.39
$ws.Range("L45").Value = 6281.1665
$ws.Range("M45").Value = -58809.39
$ws.Range("N45").Value = -7035.1665
$ws.Range("H61").Value = 3104036
$ws.Range("I61").Value = 9032.083000000001
$ws.Range("K61").Value = 9032.083000000001
$ws.Range("M61").Value = -8820.083000000001
$ws.Range("H63").Value = 16420.916
$ws.Range("I63").Value = 6817
$ws.Range("K63").Value = 6817
$ws.Range("M63").Value = -6131
$ws.Range("H66").Value = 16420.916
$ws.Range("I66").Value = 6817
$ws.Range("K66").Value = 34085
$ws.Range("M66").Value = -30653
$ws.Range("H74").Value = 830487.4399999999
$ws.Range("I74").Value = 7084.0835
$ws.Range("J74").Value = 1489210.1
$ws.Range("K74").Value = 7084.0835
$ws.Range("L74").Value = 1489210.1
$ws.Range("M74").Value = -6210.0835
$ws.Range("N74").Value = -1490958.1
$ws.Range("H77").Value = 830487.4399999999
$ws.Range("I77").Value = 7084.0835
$ws.Range("J77").Value = 1489210.1
$ws.Range("K77").Value = 35420.4175
$ws.Range("L77").Value = 7446050.5
$ws.Range("M77").Value = -31052.4175
$ws.Range("N77").Value = -7454786.5
$ws.Range("H88").Value = 2308.077
$ws.Range("I88").Value = 1301.5714
$ws.Range("K88").Value = 1301.5714
$ws.Range("M88").Value = -895.5714
$ws.Range("H91").Value = 2308.077
$ws.Range("I91").Value = 1301.5714
$ws.Range("K91").Value = 1301.5714
$ws.Range("M91").Value = 102.4286
$ws.Range("H96").Value = 39666.332
$ws.Range("J96").Value = 39666.332
$ws.Range("L96").Value = 39666.332
$ws.Range("N96").Value = -45158.332
$ws.Range("H97").Value = 9363
$ws.Range("I97").Value = 10247.454
$ws.Range("J97").Value = 4498.5
$ws.Range("K97").Value = 10247.454
$ws.Range("L97").Value = 4498.5
$ws.Range("M97").Value = -9751.454
$ws.Range("N97").Value = -5490.5
$ws.Range("H102").Value = 2304.75
$ws.Range("I102").Value = 1739.6666
$ws.Range("K102").Value = 1739.6666
$ws.Range("M102").Value = -117.6666
$ws.Range("H110").Value = 902.5
$ws.Range("I110").Value = 848.75
$ws.Range("J110").Value = 956.25
$ws.Range("K110").Value = 848.75
$ws.Range("L110").Value = 956.25
$ws.Range("M110").Value = 1196.25
$ws.Range("N110").Value = -5046.25
$ws.Range("H132").Value = 2920.2927
$ws.Range("I132").Value = 1054
$ws.Range("K132").Value = 3162
$ws.Range("M132").Value = -632
$ws.Range("H136").Value = 3104036
$ws.Range("I136").Value = 9032.083000000001
$ws.Range("K136").Value = 27096.249
$ws.Range("M136").Value = -24546.249
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 3850.6667
$ws.Range("J12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("N12").Value = -4336
$ws.Range("H86").Value = 4099.3335
$ws.Range("I86").Value = 2613.4285
$ws.Range("J86").Value = 5399.5
$ws.Range("K86").Value = 2613.4285
$ws.Range("L86").Value = 5399.5
$ws.Range("M86").Value = -1490.4285
$ws.Range("N86").Value = -7645.5
$ws.Range("H89").Value = 4099.3335
$ws.Range("I89").Value = 2613.4285
$ws.Range("J89").Value = 5399.5
$ws.Range("K89").Value = 13067.1425
$ws.Range("L89").Value = 26997.5
$ws.Range("M89").Value = -7451.1425
$ws.Range("N89").Value = -38229.5
$ws.Range("H94").Value = 1227.775
$ws.Range("I94").Value = 1133.1428
$ws.Range("K94").Value = 1133.1428
$ws.Range("M94").Value = -682.1428000000001
$ws.Range("H99").Value = 3396.4119
$ws.Range("J99").Value = 2088.6
$ws.Range("L99").Value = 2088.6
$ws.Range("N99").Value = -5084.6
$ws.Range("H105").Value = 7694.6665
$ws.Range("I105").Value = 10210.818
$ws.Range("J105").Value = 4926.9
$ws.Range("K105").Value = 10210.818
$ws.Range("L105").Value = 4926.9
$ws.Range("M105").Value = -8463.817999999999
$ws.Range("N105").Value = -8420.9
$ws.Range("H134").Value = 26473796
$ws.Range("I134").Value = 2391.9565
$ws.Range("K134").Value = 7175.869499999999
$ws.Range("M134").Value = -4640.869499999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6671.375
$ws.Range("I16").Value = 4596.75
$ws.Range("K16").Value = 4596.75
$ws.Range("M16").Value = -4309.75
$ws.Range("H31").Value = 2608.5513
$ws.Range("I31").Value = 1719.4186
$ws.Range("J31").Value = 3700.9143
$ws.Range("K31").Value = 1719.4186
$ws.Range("L31").Value = 3700.9143
$ws.Range("M31").Value = -1424.4186
$ws.Range("N31").Value = -4290.9143
$ws.Range("H34").Value = 2608.5513
$ws.Range("I34").Value = 1719.4186
$ws.Range("J34").Value = 3700.9143
$ws.Range("K34").Value = 1719.4186
$ws.Range("L34").Value = 3700.9143
$ws.Range("M34").Value = -1517.4186
$ws.Range("N34").Value = -4104.9143
$ws.Range("H58").Value = 4009.5
$ws.Range("H105").Value = 1648.2778
$ws.Range("I105").Value = 1164.9166
$ws.Range("J105").Value = 2615
$ws.Range("K105").Value = 1164.9166
$ws.Range("L105").Value = 2615
$ws.Range("M105").Value = 582.0834
$ws.Range("N105").Value = -6109
$ws.Range("H113").Value = 6671.375
$ws.Range("I113").Value = 4596.75
$ws.Range("K113").Value = 4596.75
$ws.Range("M113").Value = -2426.75
$ws.Range("H122").Value = 3621.16
$ws.Range("I122").Value = 3951.5908
$ws.Range("J122").Value = 1198
$ws.Range("K122").Value = 11854.7724
$ws.Range("L122").Value = 3594
$ws.Range("M122").Value = -9404.7724
$ws.Range("N122").Value = -8494
$ws.Range("H131").Value = 47800
$ws.Range("J131").Value = 47800
$ws.Range("L131").Value = 47800
$ws.Range("N131").Value = -57880
$ws.Range("H134").Value = 3150.9285
$ws.Range("I134").Value = 2699.2856
$ws.Range("K134").Value = 8097.8568
$ws.Range("M134").Value = -5562.8568
$ws.Range("H135").Value = 189990
$ws.Range("J135").Value = 189990
$ws.Range("L135").Value = 189990
$ws.Range("N135").Value = -200130
$ws.Range("H136").Value = 4009.5
$ws.Range("H138").Value = 102913.336
$ws.Range("J138").Value = 102913.336
$ws.Range("L138").Value = 102913.336
$ws.Range("N138").Value = -113193.336
$ws.Range("H141").Value = 427978.75
$ws.Range("I141").Value = 373332.66
$ws.Range("J141").Value = 446194.12
$ws.Range("K141").Value = 373332.66
$ws.Range("L141").Value = 446194.12
$ws.Range("M141").Value = -368152.66
$ws.Range("N141").Value = -456554.12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 439.15384
$ws.Range("I5").Value = 439.15384
$ws.Range("K5").Value = 1317.46152
$ws.Range("M5").Value = -1205.46152
$ws.Range("H23").Value = 333399.66
$ws.Range("I23").Value = 99
$ws.Range("K23").Value = 297
$ws.Range("M23").Value = -62
$ws.Range("H26").Value = 965.8
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 965.8
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 2897.4
$ws.Range("N26").Value = -3473.4
$ws.Range("M26").ClearContents()
$ws.Range("H40").Value = 3236.7778
$ws.Range("I40").Value = 5182.75
$ws.Range("K40").Value = 20731
$ws.Range("M40").Value = -20662
$ws.Range("H92").Value = 95.666664
$ws.Range("H97").Value = 365.2
$ws.Range("I97").Value = 471.85715
$ws.Range("K97").Value = 1415.57145
$ws.Range("M97").Value = -919.5714499999999
$ws.Range("H98").Value = 687
$ws.Range("J98").Value = 750.7778
$ws.Range("L98").Value = 2252.3334
$ws.Range("N98").Value = -5248.3334
$ws.Range("H104").Value = 6158.5405
$ws.Range("J104").Value = 6179.6113
$ws.Range("L104").Value = 18538.8339
$ws.Range("N104").Value = -23780.8339
$ws.Range("H106").Value = 11541.857
$ws.Range("J106").Value = 14158.6
$ws.Range("L106").Value = 42475.8
$ws.Range("N106").Value = -44367.8
$ws.Range("H109").Value = 7545.4
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 7545.4
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 22636.2
$ws.Range("N109").Value = -24716.2
$ws.Range("M109").ClearContents()
$ws.Range("H121").Value = 2121.889
$ws.Range("J121").Value = 3824.25
$ws.Range("L121").Value = 11472.75
$ws.Range("N121").Value = -14092.75
$ws.Range("H124").Value = 3301.7646
$ws.Range("I124").Value = 1815
$ws.Range("K124").Value = 5445
$ws.Range("M124").Value = -535
$ws.Range("H131").Value = 4514237
$ws.Range("J131").Value = 3205820.2
$ws.Range("L131").Value = 9617460.600000001
$ws.Range("N131").Value = -9627540.600000001
$ws.Range("H135").Value = 439.15384
$ws.Range("I135").Value = 439.15384
$ws.Range("K135").Value = 3952.38456
$ws.Range("M135").Value = -1417.38456
$ws.Range("H137").Value = 5339.4165
$ws.Range("J137").Value = 5880.533
$ws.Range("L137").Value = 17641.599
$ws.Range("N137").Value = -27841.599
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 18399.8
$ws.Range("I41").Value = 18399.8
$ws.Range("K41").Value = 18399.8
$ws.Range("M41").Value = -18044.8
$ws.Range("H57").Value = 35999.668
$ws.Range("J57").Value = 35999.668
$ws.Range("L57").Value = 35999.668
$ws.Range("N57").Value = -37639.668
$ws.Range("H80").Value = 9421529
$ws.Range("I80").Value = 169479.4
$ws.Range("J80").Value = 62620812
$ws.Range("K80").Value = 169479.4
$ws.Range("L80").Value = 62620812
$ws.Range("M80").Value = -168481.4
$ws.Range("N80").Value = -62622808
$ws.Range("H83").Value = 9421529
$ws.Range("I83").Value = 169479.4
$ws.Range("J83").Value = 62620812
$ws.Range("K83").Value = 847397
$ws.Range("L83").Value = 313104060
$ws.Range("M83").Value = -842405
$ws.Range("N83").Value = -313114044
$ws.Range("H122").Value = 11383.777
$ws.Range("I122").Value = 8649.75
$ws.Range("K122").Value = 25949.25
$ws.Range("M122").Value = -23499.25
$ws.Range("H130").Value = 132495
$ws.Range("J130").Value = 132495
$ws.Range("L130").Value = 132495
$ws.Range("N130").Value = -142535
$ws.Range("H132").Value = 11029467
$ws.Range("I132").Value = 4334.9375
$ws.Range("K132").Value = 13004.8125
$ws.Range("M132").Value = -10474.8125
$ws.Range("H141").Value = 48986.332
$ws.Range("J141").Value = 48986.332
$ws.Range("L141").Value = 48986.332
$ws.Range("N141").Value = -59346.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11299.4375
$ws.Range("J7").Value = 15294.5
$ws.Range("L7").Value = 15294.5
$ws.Range("N7").Value = -15518.5
$ws.Range("H22").Value = 6772.778
$ws.Range("J22").Value = 7369.375
$ws.Range("L22").Value = 7369.375
$ws.Range("N22").Value = -7959.375
$ws.Range("H27").Value = 6772.778
$ws.Range("J27").Value = 7369.375
$ws.Range("L27").Value = 7369.375
$ws.Range("N27").Value = -7583.375
$ws.Range("H32").Value = 3962.4
$ws.Range("I32").Value = 3962.4
$ws.Range("K32").Value = 3962.4
$ws.Range("M32").Value = -3645.4
$ws.Range("H46").Value = 8223.6875
$ws.Range("I46").Value = 50499
$ws.Range("K46").Value = 50499
$ws.Range("M46").Value = -50311
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H55").Value = 207.6875
$ws.Range("I55").Value = 165
$ws.Range("J55").Value = 217.53847
$ws.Range("K55").Value = 165
$ws.Range("L55").Value = 217.53847
$ws.Range("M55").Value = 8
$ws.Range("N55").Value = -563.53847
$ws.Range("H61").Value = 9526782
$ws.Range("I61").Value = 13336276
$ws.Range("K61").Value = 13336276
$ws.Range("M61").Value = -13336074
$ws.Range("H68").Value = 2418.3076
$ws.Range("I68").Value = 2226.7778
$ws.Range("J68").Value = 2849.25
$ws.Range("K68").Value = 2226.7778
$ws.Range("L68").Value = 2849.25
$ws.Range("M68").Value = -1477.7778
$ws.Range("N68").Value = -4347.25
$ws.Range("H71").Value = 2418.3076
$ws.Range("I71").Value = 2226.7778
$ws.Range("J71").Value = 2849.25
$ws.Range("K71").Value = 11133.889
$ws.Range("L71").Value = 14246.25
$ws.Range("M71").Value = -7389.888999999999
$ws.Range("N71").Value = -21734.25
$ws.Range("H74").Value = 44750
$ws.Range("H77").Value = 44750
$ws.Range("H93").Value = 1602.2941
$ws.Range("I93").Value = 1558
$ws.Range("K93").Value = 1558
$ws.Range("M93").Value = -310
$ws.Range("H113").Value = 9526782
$ws.Range("I113").Value = 13336276
$ws.Range("K113").Value = 13336276
$ws.Range("M113").Value = -13334106
$ws.Range("H122").Value = 4735.825
$ws.Range("I122").Value = 3620.4736
$ws.Range("K122").Value = 10861.4208
$ws.Range("M122").Value = -8411.4208
$ws.Range("H126").Value = 11299.4375
$ws.Range("J126").Value = 15294.5
$ws.Range("L126").Value = 45883.5
$ws.Range("N126").Value = -50823.5
$ws.Range("H132").Value = 8711.85
$ws.Range("I132").Value = 2302.5454
$ws.Range("K132").Value = 6907.6362
$ws.Range("M132").Value = -4377.6362
$ws.Range("H136").Value = 3567.6365
$ws.Range("J136").Value = 5685
$ws.Range("L136").Value = 17055
$ws.Range("N136").Value = -22155
$ws.Range("H140").Value = 92734.64999999999
$ws.Range("J140").Value = 92734.64999999999
$ws.Range("L140").Value = 92734.64999999999
$ws.Range("N140").Value = -103094.65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 31250
$ws.Range("I6").Value = 2500
$ws.Range("J6").Value = 60000
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 60000
$ws.Range("M6").Value = -2385
$ws.Range("N6").Value = -60230
$ws.Range("H19").Value = 8099.6
$ws.Range("I19").Value = 5124.5
$ws.Range("K19").Value = 5124.5
$ws.Range("M19").Value = -4950.5
$ws.Range("H23").Value = 7596.5
$ws.Range("J23").Value = 9996
$ws.Range("L23").Value = 9996
$ws.Range("N23").Value = -10454
$ws.Range("H70").Value = 46666.668
$ws.Range("J70").Value = 46666.668
$ws.Range("L70").Value = 46666.668
$ws.Range("N70").Value = -47296.668
$ws.Range("H73").Value = 46666.668
$ws.Range("J73").Value = 46666.668
$ws.Range("L73").Value = 46666.668
$ws.Range("N73").Value = -48850.668
$ws.Range("H104").Value = 32550.166
$ws.Range("J104").Value = 32550.166
$ws.Range("L104").Value = 32550.166
$ws.Range("N104").Value = -39538.166
$ws.Range("H123").Value = 91261.336
$ws.Range("J123").Value = 91261.336
$ws.Range("L123").Value = 91261.336
$ws.Range("N123").Value = -101061.336
$ws.Range("H132").Value = 23502.105
$ws.Range("I132").Value = 30444.086
$ws.Range("K132").Value = 91332.258
$ws.Range("M132").Value = -88802.258
$ws.Range("H140").Value = 82464.836
$ws.Range("J140").Value = 82464.836
$ws.Range("L140").Value = 82464.836
$ws.Range("N140").Value = -92824.836
$ws.Range("H141").Value = 115710.46
$ws.Range("J141").Value = 118566.91
$ws.Range("L141").Value = 118566.91
$ws.Range("N141").Value = -128926.91
